# Add two new ticket rows (68 and 69) to the worksheet, mirroring the
# existing "inline string" text format used throughout the sheet
# (dates, times and durations are stored as plain text, not as
# native Excel date/time numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds dates formatted like "2024-05-20"; Excel's COM layer
# auto-recognises that pattern and would silently convert it into a
# date serial number on assignment. Force plain text on just that
# column so the value is kept verbatim as a string, matching how every
# other row in this sheet stores its date (and the rest of the columns,
# which are not auto-detected as dates/times here).
$ws.Range("A68:A69").NumberFormat = "@"

# Row 68
$ws.Range("A68").Value = "2024-05-20"
$ws.Range("B68").Value = "10:46:30"
$ws.Range("C68").Value = "-"
$ws.Range("D68").Value = "-"
$ws.Range("E68").Value = "Etiquetadora"
$ws.Range("F68").Value = "-"
$ws.Range("G68").Value = "-"
$ws.Range("H68").Value = "10:46:37"
$ws.Range("I68").Value = "0:00:07"

# Row 69
$ws.Range("A69").Value = "2024-05-20"
$ws.Range("B69").Value = "10:47:05"
$ws.Range("C69").Value = "-"
$ws.Range("D69").Value = "-"
$ws.Range("E69").Value = "Atasco tuerca"
$ws.Range("F69").Value = "-"
$ws.Range("G69").Value = "-"
$ws.Range("H69").Value = "10:47:07"
$ws.Range("I69").Value = "0:00:02"
